# quicker data collection: runner file that loops from SE from 0 to 90
# Updates the simulation run data captured in the workbook's four data
# sheets (Radar Data, Satellite Data, Command Data, Battery Data) with the
# results of a re-run, and keeps a couple of selection/used-range tweaks
# that came along with it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Radar Data (sheet2): three side-by-side ID/Detect-Time column pairs.
# Pair A:B only has 2 data rows; pair C:D grows from 6 to 9 rows; pair
# E:F shrinks from 6 to 3 rows (rows 5-7 cleared).
# ---------------------------------------------------------------------
$wsRadar = $wb.Worksheets.Item("Radar Data")

$wsRadar.Cells.Item(2,1).Value = 3
$wsRadar.Cells.Item(2,2).Value = 45
$wsRadar.Cells.Item(2,3).Value = 1
$wsRadar.Cells.Item(2,4).Value = 13
$wsRadar.Cells.Item(2,5).Value = 8
$wsRadar.Cells.Item(2,6).Value = 85

$wsRadar.Cells.Item(3,1).Value = 7
$wsRadar.Cells.Item(3,2).Value = 83
$wsRadar.Cells.Item(3,3).Value = 2
$wsRadar.Cells.Item(3,4).Value = 24
$wsRadar.Cells.Item(3,5).Value = 10
$wsRadar.Cells.Item(3,6).Value = 105

$wsRadar.Cells.Item(4,3).Value = 3
$wsRadar.Cells.Item(4,4).Value = 37
$wsRadar.Cells.Item(4,5).Value = 9
$wsRadar.Cells.Item(4,6).Value = 109

$wsRadar.Cells.Item(5,3).Value = 4
$wsRadar.Cells.Item(5,4).Value = 49
$wsRadar.Cells.Item(5,5).ClearContents() | Out-Null
$wsRadar.Cells.Item(5,6).ClearContents() | Out-Null

$wsRadar.Cells.Item(6,3).Value = 5
$wsRadar.Cells.Item(6,4).Value = 54
$wsRadar.Cells.Item(6,5).ClearContents() | Out-Null
$wsRadar.Cells.Item(6,6).ClearContents() | Out-Null

$wsRadar.Cells.Item(7,3).Value = 6
$wsRadar.Cells.Item(7,4).Value = 67
$wsRadar.Cells.Item(7,5).ClearContents() | Out-Null
$wsRadar.Cells.Item(7,6).ClearContents() | Out-Null

$wsRadar.Cells.Item(8,3).Value = 7
$wsRadar.Cells.Item(8,4).Value = 82

$wsRadar.Cells.Item(9,3).Value = 8
$wsRadar.Cells.Item(9,4).Value = 92

$wsRadar.Cells.Item(10,3).Value = 9
$wsRadar.Cells.Item(10,4).Value = 95

# ---------------------------------------------------------------------
# Satellite Data (sheet3): Cues step table, a handful of values shift by
# one cue.
# ---------------------------------------------------------------------
$wsSat = $wb.Worksheets.Item("Satellite Data")

$wsSat.Cells.Item(14,2).Value = 1
$wsSat.Cells.Item(33,2).Value = 3
$wsSat.Cells.Item(53,2).Value = 5
$wsSat.Cells.Item(54,2).Value = 5
$wsSat.Cells.Item(63,2).Value = 6
$wsSat.Cells.Item(64,2).Value = 6
$wsSat.Cells.Item(73,2).Value = 8
$wsSat.Cells.Item(74,2).Value = 8
$wsSat.Cells.Item(83,2).Value = 8
$wsSat.Cells.Item(93,2).Value = 9

# ---------------------------------------------------------------------
# Command Data (sheet4): Assign Times re-run.
# ---------------------------------------------------------------------
$wsCmd = $wb.Worksheets.Item("Command Data")

$wsCmd.Cells.Item(2,2).Value = 14
$wsCmd.Cells.Item(3,2).Value = 25
$wsCmd.Cells.Item(4,2).Value = 38
$wsCmd.Cells.Item(5,2).Value = 51
$wsCmd.Cells.Item(8,2).Value = 83
$wsCmd.Cells.Item(9,2).Value = 86
$wsCmd.Cells.Item(10,2).Value = 96
$wsCmd.Cells.Item(11,2).Value = 106

# ---------------------------------------------------------------------
# Battery Data (sheet5): four side-by-side ID/Intercept-Time column
# pairs shrink to three (G:H cleared); C:D grows from 4 to 5 rows.
# ---------------------------------------------------------------------
$wsBat = $wb.Worksheets.Item("Battery Data")

$wsBat.Cells.Item(2,1).Value = 3
$wsBat.Cells.Item(2,2).Value = 52
$wsBat.Cells.Item(2,3).Value = 1
$wsBat.Cells.Item(2,4).Value = 22
$wsBat.Cells.Item(2,5).Value = 8
$wsBat.Cells.Item(2,6).Value = 92
$wsBat.Cells.Item(2,7).ClearContents() | Out-Null
$wsBat.Cells.Item(2,8).ClearContents() | Out-Null

$wsBat.Cells.Item(3,1).Value = 7
$wsBat.Cells.Item(3,2).Value = 91
$wsBat.Cells.Item(3,3).Value = 2
$wsBat.Cells.Item(3,4).Value = 33
$wsBat.Cells.Item(3,5).Value = 9
$wsBat.Cells.Item(3,6).Value = 109

$wsBat.Cells.Item(4,3).Value = 4
$wsBat.Cells.Item(4,4).Value = 56
$wsBat.Cells.Item(4,5).Value = 10
$wsBat.Cells.Item(4,6).Value = 110

$wsBat.Cells.Item(5,3).Value = 5
$wsBat.Cells.Item(5,4).Value = 63
$wsBat.Cells.Item(5,5).ClearContents() | Out-Null
$wsBat.Cells.Item(5,6).ClearContents() | Out-Null

$wsBat.Cells.Item(6,3).Value = 6
$wsBat.Cells.Item(6,4).Value = 78

# ---------------------------------------------------------------------
# Selection tweaks that came along with the data refresh. The workbook's
# active sheet (Command Data, tab index 3) is left unchanged at the end.
# ---------------------------------------------------------------------
$wsRadar.Activate()
$wsRadar.Range("E2:F4").Select() | Out-Null

$wsBat.Activate()
$wsBat.Range("A1:H1").Select() | Out-Null

$wsCmd.Activate()
